# Thesis progress workbook - "design" chapter figures added.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New actual-pages figures for the two days that got design-chapter work.
$ws.Range("E12").Value = 40
$ws.Range("E13").Value = 44

# "Where?" column: row 12 continues under "2.2 data visualization",
# row 13 starts the new "(Design) 2.2 data viz" entry (new shared string).
$ws.Range("H12").Value = "2.2 data visualization"
$ws.Range("H13").Value = "(Design) 2.2 data viz"

# J12 got typed over directly, breaking it out of the J9:J30 shared formula.
$ws.Range("J12").Formula = "=J11+2"

# Move the active selection to where the author left off editing.
[void]$ws.Range("H17").Select()
